$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3 was "F"; the sex value for hoge2 is corrected to "O".
$ws.Range("C3").Value = "O"

# Add a new row of data: hoge4 / 100 / T.
# Copy formatting from the row above (A4:C4) first so the new row picks up
# the same style (border/fill) already used by the other data rows, then
# fill in the values.
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)
$ws.Range("A5").Value = "hoge4"
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = "T"

# The data column (B) got a hair wider to fit the new "100" value.
$ws.Columns.Item(2).ColumnWidth = 3.5

# Land the selection back on A1 (top-left), clearing the old E4 selection.
$ws.Range("A1").Select() | Out-Null
